# Results.xlsx update
#  - Converts the "Single Color Channel Graph" chart from a line chart to a
#    clustered column (bar) chart.
#  - Adds a "Gemiddelden" (averages) summary block below the raw data:
#    per-column averages (row 58) under repeated headers (rows 55-57), plus
#    overall averages for each graph's 4 series (rows 61-63).
#  - Moves the active selection to Q14 (and drops the old frozen/scrolled
#    topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Chart: Single Color Channel Graph -> clustered column chart
# ---------------------------------------------------------------------
$chart = $ws.ChartObjects().Item(1).Chart
$chart.ChartType = 51                 # xlColumnClustered
$chart.ChartGroups().Item(1).GapWidth = 150

# ---------------------------------------------------------------------
# 2. New "Gemiddelden" (averages) block, rows 55-63
# ---------------------------------------------------------------------

# Row 55 - section title
$ws.Range("A55:B55").Merge()
$ws.Range("A55:B55").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A55").Value = "Gemiddelden"

# Row 56 - graph titles (repeat of row 1)
$ws.Range("A56:H56").Merge()
$ws.Range("A56:H56").HorizontalAlignment = -4108
$ws.Range("A56").Value = "Single Color Channel"

$ws.Range("I56:P56").Merge()
$ws.Range("I56:P56").HorizontalAlignment = -4108
$ws.Range("I56").Value = "Luminace"

# Row 57 - series titles (repeat of row 2)
$ws.Range("A57:B57").Merge()
$ws.Range("A57:B57").HorizontalAlignment = -4108
$ws.Range("A57").Value = "Female-1"

$ws.Range("C57:D57").Merge()
$ws.Range("C57:D57").HorizontalAlignment = -4108
$ws.Range("C57").Value = "Child-1"

$ws.Range("E57:F57").Merge()
$ws.Range("E57:F57").HorizontalAlignment = -4108
$ws.Range("E57").Value = "Male-2"

$ws.Range("G57:H57").Merge()
$ws.Range("G57:H57").HorizontalAlignment = -4108
$ws.Range("G57").Value = "Male-3"

$ws.Range("I57:J57").Merge()
$ws.Range("I57:J57").HorizontalAlignment = -4108
$ws.Range("I57").Value = "Female-1"

$ws.Range("K57:L57").Merge()
$ws.Range("K57:L57").HorizontalAlignment = -4108
$ws.Range("K57").Value = "Child-1"

$ws.Range("M57:N57").Merge()
$ws.Range("M57:N57").HorizontalAlignment = -4108
$ws.Range("M57").Value = "Male-2"

$ws.Range("O57:P57").Merge()
$ws.Range("O57:P57").HorizontalAlignment = -4108
$ws.Range("O57").Value = "Male-3"

# Row 58 - per-column averages
$ws.Range("A58").Formula = "=SUM(A3:A52)/50"
$ws.Range("C58").Formula = "=SUM(C3:C52)/50"
$ws.Range("E58").Formula = "=SUM(E3:E52)/50"
$ws.Range("G58").Formula = "=SUM(G3:G52)/50"
$ws.Range("I58").Formula = "=SUM(I3:I52)/50"
$ws.Range("K58").Formula = "=SUM(K3:K52)/50"
$ws.Range("M58").Formula = "=SUM(M3:M52)/50"
$ws.Range("O58").Formula = "=SUM(O3:O52)/50"

# Row 61 - "Single Color" overall-average label
$ws.Range("A61:C61").Merge()
$ws.Range("A61:C61").HorizontalAlignment = -4108
$ws.Range("A61").Value = "Gemiddelde bij single color"

# Row 62 - single color overall average + "Luminace" overall-average label
$ws.Range("A62").Formula = "=SUM(A58,C58,E58,G58)/4"

$ws.Range("I62:K62").Merge()
$ws.Range("I62:K62").HorizontalAlignment = -4108
$ws.Range("I62").Value = "Gemiddelde bij Luminace"

# Row 63 - Luminace overall average
$ws.Range("I63").Formula = "=SUM(I58,K58,M58,O58)/4"

# ---------------------------------------------------------------------
# 3. Selection
# ---------------------------------------------------------------------
$ws.Range("Q14").Select()
